# Chris_Jacobi_2026-01-19.xlsx edit
# fix: revert admin dev default; seed customers only when table empty;
#      autosave on customer select when hours/day present
#
# Concretely, this:
#  - Renames several clients across both the "Weekly Timesheet" sheet and
#    the "Jason Schema" mirror sheet.
#  - Fixes the 2026-01-21 row's hours back to 8 (was showing a stale 10).
#  - Replaces the 2026-01-25 row with a 2026-01-23 row for "Campbell".
#  - Zeroes out the rate/total columns (admin dev-default revert) and
#    updates the SUBTOTAL hours/label accordingly.
#  - Adds HOURLY SUBTOTAL / ADMIN SUBTOTAL / GRAND TOTAL rows with new
#    banding colors (and a red, bold font for the grand total row).
#  - Updates the stored employee id.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")
$ws2 = $wb.Worksheets.Item("Jason Schema")

# ---------------------------------------------------------------------
# Weekly Timesheet (sheet1) — per-row client renames + rate/total reset
# ---------------------------------------------------------------------

# Row 2: 2026-01-19
$ws1.Range("B2").Value = "Smithers"
$ws1.Range("E2").Value = 0
$ws1.Range("F2").Value = 0

# Row 3: 2026-01-20
$ws1.Range("B3").Value = "Bottomley"
$ws1.Range("E3").Value = 0
$ws1.Range("F3").Value = 0

# Row 4: 2026-01-21
$ws1.Range("B4").Value = "Behrens"
$ws1.Range("C4").Value = 8
$ws1.Range("E4").Value = 0
$ws1.Range("F4").Value = 0

# Row 5: 2026-01-22
$ws1.Range("B5").Value = "Goodrich"
$ws1.Range("E5").Value = 0
$ws1.Range("F5").Value = 0

# Row 6: was 2026-01-25 / McGill -> now 2026-01-23 / Campbell
# (force text entry so Excel doesn't auto-coerce the date-shaped string
#  into a date serial, then drop back to the Normal style so no stray
#  number format sticks to the cell)
$ws1.Range("A6").NumberFormat = "@"
$ws1.Range("A6").Value = "2026-01-23"
$ws1.Range("A6").Style = "Normal"
$ws1.Range("B6").Value = "Campbell"
$ws1.Range("E6").Value = 0
$ws1.Range("F6").Value = 0

# SUBTOTAL row
$ws1.Range("C8").Value = 40
$ws1.Range("D8").Value = "Reg: 40 / OT: 0"
$ws1.Range("F8").Value = 0

# ---------------------------------------------------------------------
# New summary rows 11-13 on the Weekly Timesheet sheet
# ---------------------------------------------------------------------

# Row 11: HOURLY SUBTOTAL - bold font on a light tan band
$hourly = $ws1.Range("A11:F11")
$hourly.Font.Bold = $true
$hourly.Interior.Color = 14742522  # RGB(FA,F3,E0) as BGR-packed OLE color
$ws1.Range("A11").Value = "HOURLY SUBTOTAL"
$ws1.Range("F11").Value = 0

# Row 12: ADMIN SUBTOTAL - same banding as row 11
$adminSub = $ws1.Range("A12:F12")
$adminSub.Font.Bold = $true
$adminSub.Interior.Color = 14742522
$ws1.Range("A12").Value = "ADMIN SUBTOTAL"
$ws1.Range("F12").Value = 0

# Row 13: GRAND TOTAL - bold red font on a light green band
$grand = $ws1.Range("A13:F13")
$grand.Font.Bold = $true
$grand.Font.Color = 255            # RGB(FF,00,00) as BGR-packed OLE color
$grand.Interior.Color = 14743784   # RGB(E8,F8,E0) as BGR-packed OLE color
$ws1.Range("A13").Value = "GRAND TOTAL"
$ws1.Range("F13").Value = 0

# ---------------------------------------------------------------------
# Jason Schema (sheet2) mirror of the same data
# ---------------------------------------------------------------------

# Row 2: 2026-01-19
$ws2.Range("D2").Value = "Smithers"
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 0

# Row 3: 2026-01-20
$ws2.Range("D3").Value = "Bottomley"
$ws2.Range("F3").Value = 0
$ws2.Range("G3").Value = 0

# Row 4: 2026-01-21
$ws2.Range("D4").Value = "Behrens"
$ws2.Range("E4").Value = 8
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Value = 0

# Row 5: 2026-01-22
$ws2.Range("D5").Value = "Goodrich"
$ws2.Range("F5").Value = 0
$ws2.Range("G5").Value = 0

# Row 6: was 2026-01-25 / McGill -> now 2026-01-23 / Campbell
$ws2.Range("C6").NumberFormat = "@"
$ws2.Range("C6").Value = "2026-01-23"
$ws2.Range("C6").Style = "Normal"
$ws2.Range("D6").Value = "Campbell"
$ws2.Range("F6").Value = 0
$ws2.Range("G6").Value = 0

# ---------------------------------------------------------------------
# Employee id update (shared across all Jason Schema rows)
# ---------------------------------------------------------------------
$ws2.Range("B2:B6").Value = "emp_lf0u97k0"
